# Update the cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new cell value, built from the upstream diff.
$changes = @{
  2  = @{ "D"="67.914.26";  "E"="  +0.17%  " }
  3  = @{ "D"="3.333.46";   "E"="  +0.45%  " }
  4  = @{ "D"="0.999";      "E"="  +0.02%  " }
  5  = @{ "D"="584.71";     "E"="  +0.48%  " }
  6  = @{ "D"="177.63";     "E"="  +1.95%  " }
  7  = @{ "E"="  +0.02%  " }
  8  = @{ "E"="  +1.72%  " }
  9  = @{ "E"="  +4.44%  " }
  10 = @{ "E"="  +1.43%  " }
  11 = @{ "D"="48.13";      "E"="  +6.28%  " }
  12 = @{ "E"="  +1.86%  " }
  13 = @{ "D"="699.03";     "E"="  +4.85%  " }
  14 = @{ "D"="3.878.93";   "E"="  +0.57%  " }
  15 = @{ "E"="  +0.96%  " }
  16 = @{ "D"="67.965.91";  "E"="  +0.08%  " }
  17 = @{ "E"="  +1.13%  " }
  18 = @{ "D"="3.339.32";   "E"="  +0.73%  " }
  19 = @{ "D"="17.52";      "E"="  +0.56%  " }
  20 = @{ "E"="  +2.77%  " }
  21 = @{ "D"="0.895";      "E"="  +0.97%  " }
  22 = @{ "D"="5.41";       "E"="  +1.04%  " }
  23 = @{ "D"="16.91";      "E"="  -0.14%  " }
  24 = @{ "D"="100.24";     "E"="  +2.67%  " }
  25 = @{ "D"="3.91";       "E"="  +1.97%  " }
  26 = @{ "E"="  +0.63%  " }
  27 = @{ "D"="9.47";       "E"="  +2.38%  " }
  28 = @{ "D"="33.03";      "E"="  -1.65%  " }
  29 = @{ "D"="8.55";       "E"="  +1.82%  " }
  30 = @{ "E"="  -4.61%  " }
  31 = @{ "D"="577.77";     "E"="  -0.97%  " }
  32 = @{ "E"="  +0.92%  " }
  33 = @{ "E"="  +1.93%  " }
  34 = @{ "D"="3.738.86";   "E"="  +0.51%  " }
  35 = @{ "D"="57.36";      "E"="  +0.32%  " }
  36 = @{ "E"="  +0.08%  " }
  37 = @{ "D"="3.37";       "E"="  +1.58%  " }
  38 = @{ "D"="35.32";      "E"="  +8.82%  " }
  39 = @{ "E"="  +2.68%  " }
  40 = @{ "B"="Fetch.AI"; "C"="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; "D"="2.63"; "E"="  +0.53%  " }
  41 = @{ "B"="Stacks";   "C"="https://coinranking.com/coin/mMPrMcB7+stacks-stx";        "D"="3.16"; "E"="  +2.58%  " }
  42 = @{ "D"="0.0₃0674"; "E"="  +1.87%  " }
  44 = @{ "D"="3.27";       "E"="  +0.14%  " }
  45 = @{ "E"="  +0.81%  " }
  46 = @{ "D"="2.62";       "E"="  +1.57%  " }
  47 = @{ "E"="  +1.28%  " }
  48 = @{ "E"="  -0.07%  " }
  49 = @{ "E"="  -1.25%  " }
  50 = @{ "D"="130.97";     "E"="  +2.59%  " }
  51 = @{ "E"="  +1.25%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $newValue = $cols[$col]
        $range = $ws.Range($cellRef)

        # All cells in this sheet (columns B-E) are stored as text in the
        # source workbook, even when the text happens to look like a plain
        # number (e.g. "584.71"). Excel's COM Value setter auto-detects
        # numeric-looking strings and stores them as numbers, so for those
        # values we temporarily force a text number format, assign the
        # value, then restore the cell's original style/format so the cell
        # keeps looking exactly as it did before (no visible style change).
        $looksNumeric = $newValue -match '^-?\d+(\.\d+)?$'

        if ($looksNumeric) {
            $origStyle = $range.Style
            $range.NumberFormat = "@"
            $range.Value = $newValue
            $range.Style = $origStyle
        } else {
            $range.Value = $newValue
        }
    }
}
